$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.324.66"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").Value = "1.935.18"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7492"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.09%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "28.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3186"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07150"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7821"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08037"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.85%  "

$ws.Range("D13").Value = "1.921.63"
$ws.Range("E13").Value = "  -0.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.413"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.60%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.95%  "

$ws.Range("D17").Value = "30.336.98"
$ws.Range("E17").Value = "  +0.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.088"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007973"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.06%  "

$ws.Range("D21").Value = "2.171.19"
$ws.Range("E21").Value = "  -0.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.676"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.570"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.38%  "

$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1303"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.199"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.370"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.548"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.427"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.155"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05292"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.334"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7592"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.788"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01954"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("E39").Value = "  +0.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.515"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4528"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.983"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8407"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.43%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.706"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.79%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.998"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.16%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "965.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.69%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1223"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.52%  "

